# Generate Report for Handoff
#
# Marks the pending handoff as complete: status moves from "In Translation"
# to "Ready for handoff" and the handoff timestamps are bumped to the new
# generation time. Column widths for the affected "Status"-type columns grow
# to fit the longer "Ready for handoff" text (was auto-fit to the old text).

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-09-04 11:04:25"

# Columns E (zh-cn) and F (de-de) widen to fit the new, longer status text.
# (16.3333... is the ColumnWidth input that resolves to the same stored
# column width - ~17.22 characters - that Excel's own auto-fit produced.)
$ws1.Columns.Item(5).ColumnWidth = 16.333333333333336
$ws1.Columns.Item(6).ColumnWidth = 16.333333333333336

# --- zh-cn sheet -----------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("H2").Value = "2016-09-04 11:04:21"

# Status column widens to fit the new, longer status text.
$ws2.Columns.Item(3).ColumnWidth = 16.333333333333336

# --- de-de sheet -----------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("H2").Value = "2016-09-04 11:04:25"

# Status column widens to fit the new, longer status text.
$ws3.Columns.Item(3).ColumnWidth = 16.333333333333336
